# A new weekly price record was added to the "Vega Modelo de Temuco -
# Espinaca" sheet. In the source data the new record sorts to the top of
# the existing date-ordered block (rows 270-315), so the whole block is
# pushed down by one row (270-315 -> 271-316) and the new record is
# written into the vacated row 270. This also grows the sheet's used
# range from A1:R315 to A1:R316.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 270:315 down to 271:316, leaving row 270 empty.
$ws.Rows(270).Insert()

# Populate the newly inserted row 270 with the new record.
$ws.Range("A270").Value = 10
$ws.Range("B270").Value = "Vega Modelo de Temuco"
$ws.Range("C270").Value = "La Araucanía"
$ws.Range("D270").Value = 45218
$ws.Range("E270").Value = 9
$ws.Range("F270").Value = 100112012
$ws.Range("G270").Value = "Espinaca"
$ws.Range("H270").Value = "Sin especificar"
$ws.Range("I270").Value = "Primera"
$ws.Range("J270").Value = 55
$ws.Range("K270").Value = 14000
$ws.Range("L270").Value = 14000
$ws.Range("M270").Value = 14000
$ws.Range("N270").Value = "$/docena de atados"
$ws.Range("O270").Value = "Región de La Araucanía"
$ws.Range("P270").Value = 4667
$ws.Range("Q270").Value = 3
$ws.Range("R270").Value = "Hortaliza"
